$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster entries for Isaiah Hartenstein (row 14) and Damian Lillard
# (row 15) move up to the top of the list, right after the header, in
# swapped order (Damian Lillard first, then Isaiah Hartenstein). Every
# other player keeps their existing relative order and simply shifts down
# by two rows; the last three rows (Kentavious Caldwell-Pope, Brandon
# Ingram, LaMelo Ball) stay exactly where they were.

# Capture the two rows that are moving before anything is overwritten.
$damianLillard = @(
    $ws.Cells.Item(15, 1).Value2,
    $ws.Cells.Item(15, 2).Value2,
    $ws.Cells.Item(15, 3).Value2
)
$isaiahHartenstein = @(
    $ws.Cells.Item(14, 1).Value2,
    $ws.Cells.Item(14, 2).Value2,
    $ws.Cells.Item(14, 3).Value2
)

# Capture the block of rows (old rows 2-13) that will shift down by two.
$shiftingBlock = @()
for ($r = 2; $r -le 13; $r++) {
    $shiftingBlock += ,@(
        $ws.Cells.Item($r, 1).Value2,
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2
    )
}

# Write Damian Lillard then Isaiah Hartenstein into the new top rows (2-3).
$ws.Cells.Item(2, 1).Value2 = $damianLillard[0]
$ws.Cells.Item(2, 2).Value2 = $damianLillard[1]
$ws.Cells.Item(2, 3).Value2 = $damianLillard[2]

$ws.Cells.Item(3, 1).Value2 = $isaiahHartenstein[0]
$ws.Cells.Item(3, 2).Value2 = $isaiahHartenstein[1]
$ws.Cells.Item(3, 3).Value2 = $isaiahHartenstein[2]

# Write the shifted block into its new home, rows 4-15.
for ($i = 0; $i -lt $shiftingBlock.Count; $i++) {
    $destRow = 4 + $i
    $ws.Cells.Item($destRow, 1).Value2 = $shiftingBlock[$i][0]
    $ws.Cells.Item($destRow, 2).Value2 = $shiftingBlock[$i][1]
    $ws.Cells.Item($destRow, 3).Value2 = $shiftingBlock[$i][2]
}

# Rows 16-18 (Kentavious Caldwell-Pope, Brandon Ingram, LaMelo Ball) are
# untouched - they already hold their correct final values.
